# Automatische test-sync: 2025-08-03 18:31:50
# Appends a new test-mail log row to the "Logs" sheet and updates the
# "Dashboard" summary sheet to reflect the new category counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 37 with the new test mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 37
$logs.Cells.Item($newRow, 1).Value = "Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #9: Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nHelaas kunnen we uw vraag niet beantwoorden omdat deze e-mail als een testmail is gemarkeerd. Mocht u een specifieke vraag hebben waarover wij u kunnen helpen, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 18:31:22"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Writing the multi-line "Antwoord" text auto-pins a custom row height;
# AutoFit() clears that flag again so the row matches the default sizing
# used by every other row in the table.
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional-formatting ranges that cover the log table so
# the newly added row 37 is included too (D/G/H/I/J columns).
$ccols = "D", "G", "H", "I", "J"
foreach ($col in $ccols) {
    $oldRange = $logs.Range($col + "2:" + $col + "36")
    $newRange = $logs.Range($col + "2:" + $col + "37")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: "Overig" now has 10 entries, "Planning / Afspraak"
#    still has 9 - swap the row order so the list stays sorted by count.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(2, 1).Value = "Overig"
$dash.Cells.Item(2, 2).Value = 10
$dash.Cells.Item(3, 1).Value = "Planning / Afspraak"
$dash.Cells.Item(3, 2).Value = 9
